# "cleanups through 2008 and 2004 Refresh"
#
# - Rename the worksheet tab to match the new slug-style name.
# - Clean up header-row label casing/wording (row 2, before the title
#   row above it gets removed).
# - Remove the leading title row ("Criminal Offenses - Noncampus") so the
#   column headers become row 1 and the data shifts up by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet.
$ws.Name = "criminal-offenses-noncampus-vir"

# Clean up the header labels (still on row 2 at this point, since the
# title row above it hasn't been removed yet).
$ws.Cells.Item(2, 1).Value  = "Survey Year"
$ws.Cells.Item(2, 2).Value  = "UnitID"
$ws.Cells.Item(2, 3).Value  = "Institution Name"
$ws.Cells.Item(2, 7).Value  = "Murder/Non-Negligent Manslaughter"
$ws.Cells.Item(2, 8).Value  = "Negligent Manslaughter"
$ws.Cells.Item(2, 9).Value  = "Sex Offenses - Forcible"
$ws.Cells.Item(2, 10).Value = "Sex Offenses - Non-Forcible"
$ws.Cells.Item(2, 12).Value = "Aggravated Assault"
$ws.Cells.Item(2, 14).Value = "Motor Vehicle Theft"

# Drop the standalone title row so the header row (now cleaned up)
# becomes row 1 and every data row shifts up by one.
$ws.Rows(1).Delete()
